# Daily attendance processing - 2026-01-01 11:01:52
# Re-sort the comma-separated "Recorded By" values (column G) alphabetically
# (case-insensitive) for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ",\s*"
        $trimmed = $parts | ForEach-Object { $_.Trim() }
        $sorted = $trimmed | Sort-Object { $_.ToLower() }
        $newValue = [string]::Join(", ", $sorted)

        if ($newValue -ne $value) {
            $cell.Value2 = $newValue
        }
    }
}
